$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on Price/Volume columns where the new value
# would otherwise be auto-parsed by Excel as a number (losing formatting
# like a trailing zero, or multi-dot "thousand.thousand.decimal" strings).

$ws.Range('D2').Value = '29.556.67'
$ws.Range('E2').Value = '  -1.38%  '
$ws.Range('D3').Value = '2.009.33'
$ws.Range('E3').Value = '  -5.05%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.010'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '331.33'
$ws.Range('E5').Value = '  -4.54%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.010'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4946'
$ws.Range('E7').Value = '  -4.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4198'
$ws.Range('E8').Value = '  -5.53%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '53.55'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08833'
$ws.Range('E10').Value = '  -5.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.120'
$ws.Range('E11').Value = '  -5.28%  '
$ws.Range('D12').Value = '2.156.16'
$ws.Range('E12').Value = '  +2.00%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.08'
$ws.Range('E13').Value = '  -8.24%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.114'
$ws.Range('E14').Value = '  -4.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.508'
$ws.Range('E15').Value = '  -5.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.36'
$ws.Range('E16').Value = '  -6.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.010'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001107'
$ws.Range('E18').Value = '  -4.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06623'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.57'
$ws.Range('E20').Value = '  -8.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.011'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.005'
$ws.Range('E22').Value = '  -4.68%  '
$ws.Range('D23').Value = '29.626.95'
$ws.Range('E23').Value = '  -1.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.88'
$ws.Range('E24').Value = '  -6.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.280'
$ws.Range('E25').Value = '  -1.74%  '
$ws.Range('D26').Value = '2.350.16'
$ws.Range('E26').Value = '  -0.42%  '
$ws.Range('E27').Value = '  -2.85%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.613'
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '20.67'
$ws.Range('E29').Value = '  -6.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.342'
$ws.Range('E30').Value = '  -7.91%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.07'
$ws.Range('E31').Value = '  -5.33%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.057'
$ws.Range('E32').Value = '  -7.73%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09970'
$ws.Range('E33').Value = '  -5.59%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.558'
$ws.Range('E34').Value = '  -12.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.811'
$ws.Range('E35').Value = '  -4.18%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.808'
$ws.Range('E36').Value = '  -7.07%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.643'
$ws.Range('E37').Value = '  -10.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02456'
$ws.Range('E38').Value = '  -6.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06391'
$ws.Range('E39').Value = '  -6.69%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.293'
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.83'
$ws.Range('E41').Value = '  -6.86%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6521'
$ws.Range('E42').Value = '  -7.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.2073'
$ws.Range('E43').Value = '  -7.49%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.010'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6336'
$ws.Range('E45').Value = '  -7.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.220'
$ws.Range('E46').Value = '  -5.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.49'
$ws.Range('E47').Value = '  -6.91%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.266'
$ws.Range('E48').Value = '  +0.51%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.564'
$ws.Range('E49').Value = '  -1.91%  '
$ws.Range('E50').Value = '  -7.78%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07030'
$ws.Range('E51').Value = '  -1.23%  '
